$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking target cells as Text so Excel keeps the literal string
# (matches source data semantics: these columns hold formatted display strings, not real numbers)
$numericLookingCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D19", "E19", "D20", "E20", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values cell by cell
$ws.Range("D2").Value = '294.52'
$ws.Range("E2").Value = '-4.82%'
$ws.Range("D3").Value = '40.09'
$ws.Range("E3").Value = '-2.64%'
$ws.Range("D4").Value = '5.024'
$ws.Range("E4").Value = '-3.76%'
$ws.Range("D5").Value = '0.07379'
$ws.Range("E5").Value = '-3.99%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '4.313'
$ws.Range("E6").Value = '-0.09%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '1.537'
$ws.Range("E7").Value = '-6.52%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9224'
$ws.Range("E8").Value = '0.84%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1179'
$ws.Range("E9").Value = '-5.28%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1771'
$ws.Range("E10").Value = '-2.97%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.08657'
$ws.Range("E11").Value = '-5.65%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.04166'
$ws.Range("E12").Value = '-1.47%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.1055'
$ws.Range("E13").Value = '0.21%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001270'
$ws.Range("E14").Value = '0.77%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005779'
$ws.Range("E15").Value = '-1.32%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.373'
$ws.Range("E16").Value = '0.69%'
$ws.Range("D19").Value = '7.603'
$ws.Range("E19").Value = '3.04%'
$ws.Range("D20").Value = '0.1355'
$ws.Range("E20").Value = '-3.30%'
$ws.Range("E21").Value = '-0.35%'
$ws.Range("E22").Value = '-4.89%'
$ws.Range("D23").Value = '0.001284'
$ws.Range("E23").Value = '1.36%'
$ws.Range("D24").Value = '0.003899'
$ws.Range("E24").Value = '-4.68%'
$ws.Range("D25").Value = '0.0001295'
$ws.Range("E25").Value = '-0.51%'
$ws.Range("D26").Value = '0.0003737'
$ws.Range("E26").Value = '-95.02%'
$ws.Range("D38").Value = '0.02308'
$ws.Range("E38").Value = '-9.66%'
$ws.Range("D39").Value = '0.04993'
$ws.Range("E39").Value = '-6.63%'
$ws.Range("D40").Value = '0.007714'
$ws.Range("E40").Value = '-1.46%'
$ws.Range("D41").Value = '0.1274'
$ws.Range("E41").Value = '-3.02%'
$ws.Range("E42").Value = '115.72%'
$ws.Range("D43").Value = '0.007400'
$ws.Range("E43").Value = '10.85%'
$ws.Range("E44").Value = '-4.57%'
$ws.Range("D45").Value = '0.3195'
$ws.Range("E45").Value = '4.07%'
$ws.Range("D46").Value = '0.00006437'
$ws.Range("E46").Value = '-4.76%'
$ws.Range("E47").Value = '0.29%'
$ws.Range("E48").Value = '13.70%'
$ws.Range("D49").Value = '0.004216'
$ws.Range("E49").Value = '35.85%'
$ws.Range("D50").Value = '0.00002108'
$ws.Range("E50").Value = '0.29%'
$ws.Range("D51").Value = '0.0002008'
$ws.Range("E51").Value = '0.29%'
